$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Alternativa 1" alternative flow: step 5 -> step 6 renumbering
$ws.Range("B20").Value = "Alternativa 1 [Peças Incompativeis] (passo 6)"
$ws.Range("D20").Value = "6.1 Verifica que especificações não estão corretas"
$ws.Range("D21").Value = "6.2 Informa que escolheu peças incompativeis e/ou peças em falta"

# The "Regressa a" reference in the second alternative flow now points back to step 1
$ws.Range("D24").Value = "Regressa a 1"

# Update the view so the selection matches the author's last position
$ws.Range("D24").Select()
